$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.806.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.13%  '

$ws.Range("D3").Value = "'2.278.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.44%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").Value = "'315.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.48%  '

$ws.Range("D6").Value = "'102.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.68%  '

$ws.Range("D7").Value = "'0.624"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.44%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").Value = "'0.601"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.29%  '

$ws.Range("D10").Value = "'38.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.58%  '

$ws.Range("D11").Value = "'0.0904"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.47%  '

$ws.Range("D12").Value = "'8.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.91%  '

$ws.Range("E13").Value = '  -0.41%  '

$ws.Range("D14").Value = "'0.952"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.32%  '

$ws.Range("D15").Value = "'15.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.82%  '

$ws.Range("D16").Value = "'2.620.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.46%  '

$ws.Range("D17").Value = "'2.273.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -10.97%  '

$ws.Range("D18").Value = "'41.761.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.17%  '

$ws.Range("D19").Value = "'7.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.77%  '

$ws.Range("E20").Value = '  -1.54%  '

$ws.Range("D21").Value = "'73.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.72%  '

$ws.Range("D22").Value = "'280.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.86%  '

$ws.Range("E23").Value = '  -4.91%  '

$ws.Range("D24").Value = "'9.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.14%  '

$ws.Range("E25").Value = '  -3.88%  '

$ws.Range("E26").Value = '  +0.74%  '

$ws.Range("D27").Value = "'10.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.43%  '

$ws.Range("D28").Value = "'2.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.02%  '

$ws.Range("D29").Value = "'22.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.44%  '

$ws.Range("D30").Value = "'163.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.95%  '

$ws.Range("D31").Value = "'34.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.13%  '

$ws.Range("D32").Value = "'0.0868"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.88%  '

$ws.Range("E33").Value = '  -0.19%  '

$ws.Range("D34").Value = "'5.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.90%  '

$ws.Range("E35").Value = '  -0.16%  '

$ws.Range("E36").Value = '  -9.54%  '

$ws.Range("D37").Value = "'4.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.18%  '

$ws.Range("D38").Value = "'2.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.49%  '

$ws.Range("D39").Value = "'0.0345"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.42%  '

$ws.Range("D40").Value = "'3.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.55%  '

$ws.Range("D41").Value = "'101.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +16.07%  '

$ws.Range("E42").Value = '  -2.98%  '

$ws.Range("D43").Value = "'68.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.84%  '

$ws.Range("E44").Value = '  +0.21%  '

$ws.Range("D45").Value = "'0.223"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.16%  '

$ws.Range("D46").Value = "'115.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.70%  '

$ws.Range("D47").Value = "'11.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.06%  '

$ws.Range("D48").Value = "'9.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.49%  '

$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").Value = "'5.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.05%  '

$ws.Range("B50").Value = 'ordi'
$ws.Range("C50").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D50").Value = "'75.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.93%  '

$ws.Range("E51").Value = '  -3.12%  '
